# Add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right before the current "2022-Q2"
#    sheet (i.e. right after "总计"), and populate it with the fund
#    holdings table for that quarter.
# 2) Update the "总计" (summary) sheet: shift the existing quarterly rows
#    down by one and insert the new 2022-Q3 summary row at the top of the
#    data (row 2).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" sheet before "2022-Q2"
# ---------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($sheetQ2)
$newSheet.Name = "2022-Q3"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "000031"
$newSheet.Range("C2").Value = "华夏复兴混合A"
$newSheet.Range("D2").Value = "24.25"
$newSheet.Range("E2").Value = "88.23"
$newSheet.Range("F2").Value = "5.23"
$newSheet.Range("G2").Value = "1.2683"
$newSheet.Range("H2").Value = 6

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "011140"
$newSheet.Range("C3").Value = "广发聚鸿六个月持有期混合E"
$newSheet.Range("D3").Value = "11.46"
$newSheet.Range("E3").Value = "93.24"
$newSheet.Range("F3").Value = "2.91"
$newSheet.Range("G3").Value = "0.3335"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "011138"
$newSheet.Range("C4").Value = "广发聚鸿六个月持有期混合A"
$newSheet.Range("D4").Value = "10.81"
$newSheet.Range("E4").Value = "93.24"
$newSheet.Range("F4").Value = "2.91"
$newSheet.Range("G4").Value = "0.3146"
$newSheet.Range("H4").Value = 10

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "162720"
$newSheet.Range("C5").Value = "广发创业板两年定期开放混合"
$newSheet.Range("D5").Value = "6.33"
$newSheet.Range("E5").Value = "94.20"
$newSheet.Range("F5").Value = "4.09"
$newSheet.Range("G5").Value = "0.2589"
$newSheet.Range("H5").Value = 6

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "501070"
$newSheet.Range("C6").Value = "广发睿阳三年定期开放混合"
$newSheet.Range("D6").Value = "6.62"
$newSheet.Range("E6").Value = "51.01"
$newSheet.Range("F6").Value = "3.85"
$newSheet.Range("G6").Value = "0.2549"
$newSheet.Range("H6").Value = 5

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "160211"
$newSheet.Range("C7").Value = "国泰中小盘成长混合（LOF）"
$newSheet.Range("D7").Value = "6.30"
$newSheet.Range("E7").Value = "87.18"
$newSheet.Range("F7").Value = "3.04"
$newSheet.Range("G7").Value = "0.1915"
$newSheet.Range("H7").Value = 7

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "161838"
$newSheet.Range("C8").Value = "银华创业板两年定期开放混合"
$newSheet.Range("D8").Value = "4.45"
$newSheet.Range("E8").Value = "96.33"
$newSheet.Range("F8").Value = "4.23"
$newSheet.Range("G8").Value = "0.1882"
$newSheet.Range("H8").Value = 10

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "320011"
$newSheet.Range("C9").Value = "诺安中小盘精选混合"
$newSheet.Range("D9").Value = "3.02"
$newSheet.Range("E9").Value = "73.05"
$newSheet.Range("F9").Value = "3.71"
$newSheet.Range("G9").Value = "0.1120"
$newSheet.Range("H9").Value = 9

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "519673"
$newSheet.Range("C10").Value = "银河康乐股票A"
$newSheet.Range("D10").Value = "1.94"
$newSheet.Range("E10").Value = "93.82"
$newSheet.Range("F10").Value = "5.00"
$newSheet.Range("G10").Value = "0.0970"
$newSheet.Range("H10").Value = 6

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "013936"
$newSheet.Range("C11").Value = "广发睿升混合A"
$newSheet.Range("D11").Value = "2.92"
$newSheet.Range("E11").Value = "88.08"
$newSheet.Range("F11").Value = "3.32"
$newSheet.Range("G11").Value = "0.0969"
$newSheet.Range("H11").Value = 9

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "320015"
$newSheet.Range("C12").Value = "诺安行业轮动混合"
$newSheet.Range("D12").Value = "1.11"
$newSheet.Range("E12").Value = "73.94"
$newSheet.Range("F12").Value = "3.87"
$newSheet.Range("G12").Value = "0.0430"
$newSheet.Range("H12").Value = 7

$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "006429"
$newSheet.Range("C13").Value = "诺安恒鑫混合"
$newSheet.Range("D13").Value = "0.64"
$newSheet.Range("E13").Value = "70.94"
$newSheet.Range("F13").Value = "3.69"
$newSheet.Range("G13").Value = "0.0236"
$newSheet.Range("H13").Value = 10

$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "016018"
$newSheet.Range("C14").Value = "银河康乐股票C"
$newSheet.Range("D14").Value = "0.35"
$newSheet.Range("E14").Value = "93.82"
$newSheet.Range("F14").Value = "5.00"
$newSheet.Range("G14").Value = "0.0175"
$newSheet.Range("H14").Value = 6

$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "011139"
$newSheet.Range("C15").Value = "广发聚鸿六个月持有期混合C"
$newSheet.Range("D15").Value = "0.56"
$newSheet.Range("E15").Value = "93.24"
$newSheet.Range("F15").Value = "2.91"
$newSheet.Range("G15").Value = "0.0163"
$newSheet.Range("H15").Value = 10

$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "013937"
$newSheet.Range("C16").Value = "广发睿升混合C"
$newSheet.Range("D16").Value = "0.28"
$newSheet.Range("E16").Value = "88.08"
$newSheet.Range("F16").Value = "3.32"
$newSheet.Range("G16").Value = "0.0093"
$newSheet.Range("H16").Value = 9

$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "015073"
$newSheet.Range("C17").Value = "华夏复兴混合C"
$newSheet.Range("D17").Value = "0.04"
$newSheet.Range("E17").Value = "88.23"
$newSheet.Range("F17").Value = "5.23"
$newSheet.Range("G17").Value = "0.0021"
$newSheet.Range("H17").Value = 6

# header row + A-column cells use the bold/boxed style (same "s=2" look as
# every other quarter sheet) - copy it over from the sheet we inserted
# next to so the new sheet matches the existing visual style exactly.
$sheetQ2.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$sheetQ2.Range("A2:A17").Copy($newSheet.Range("A2:A17"))
$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 2: update "总计" - push existing rows down one and write the new
# 2022-Q3 summary row at the top
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Range("A2:D8").Copy($totalSheet.Range("A3:D9"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 16
$totalSheet.Range("D2").Value = 3.23

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("A8").Value = 6
$totalSheet.Range("A9").Value = 7

$totalSheet.Range("A1").Select()
